$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rename the table headers (English instead of Finnish):
#   Osa                  -> Part
#   Ajoitus (viikko)      -> Deadline
#   Valmistunut viikolla  -> Finished
#   Tehtäviä tehty        -> Completed tasks
# Writing straight into the header cells of the ListObject keeps the table
# column names in sync automatically.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "Part"
$ws.Range("C2").Value = "Deadline"
$ws.Range("D2").Value = "Finished"
$ws.Range("E2").Value = "Completed tasks"

# ---------------------------------------------------------------------------
# Replace the "week number" schedule in column C with real deadline dates,
# mark part 0 as finished (column D) and record that 6 tasks were completed
# for it (column E).
#
# A handful of the dd.mm.yyyy strings below are ambiguous (both day and
# month <= 12), so a plain .Value assignment would get silently parsed into
# a date serial number. Force those particular cells to Text first, then
# copy the *formatting only* back from an untouched sibling cell so the
# cell keeps using the original (shared) style instead of a brand-new one.
# ---------------------------------------------------------------------------
function Set-LiteralDate($addr, $text) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
}

# Unambiguous dates - a normal value assignment is never reinterpreted.
$ws.Range("C4").Value = "19.06.2022"
$ws.Range("C7").Value = "17.7.2022"
$ws.Range("C8").Value = "24.7.2022"
$ws.Range("C9").Value = "31.7.2022"
$ws.Range("C11").Value = "14.8.2022"
$ws.Range("C12").Value = "21.8.2022"
$ws.Range("C13").Value = "28.8.2022"
$ws.Range("C16").Value = "18.9.2022"
$ws.Range("D3").Value = "13.06.2022"

# Ambiguous dates - force text, fix the style up afterwards.
Set-LiteralDate "C3"  "12.06.2022"
Set-LiteralDate "C5"  "3.7.2022"
Set-LiteralDate "C6"  "10.7.2022"
Set-LiteralDate "C10" "7.8.2022"
Set-LiteralDate "C14" "4.9.2022"
Set-LiteralDate "C15" "11.9.2022"

$ws.Range("C7").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Part 0 is now completed, with 6 tasks done.
$ws.Range("E3").Value = 6

# ---------------------------------------------------------------------------
# The old "continue part 0's tasks" reminder note next to the table is no
# longer needed.
# ---------------------------------------------------------------------------
$ws.Range("G3").ClearContents()

# ---------------------------------------------------------------------------
# Point the totals-row formula at the renamed "Completed tasks" column.
# ---------------------------------------------------------------------------
$ws.Range("E17").Formula = "=SUM(Table1[Completed tasks])"

# Column E needs to be a little wider for the new "Completed tasks" header.
$ws.Columns("E").ColumnWidth = 17.35

# Move the selection like in the saved workbook.
$ws.Range("G10").Select()
